# Update the BOM generation date on every sheet (All, Top, Bottom)
# from "2022 September 15" to "2022 September 28".
$wb = $excel.ActiveWorkbook

foreach ($ws in $wb.Worksheets) {
    $cell = $ws.Range("A2")
    if ($cell.Value2 -eq "2022 September 15") {
        $cell.Value = "2022 September 28"
    }
}
